$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-03"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 06-03)"

# Update June total (row 7) in the "Total"/2022 column
$ws.Range("I7").Value = 9

# Update the grand Total row (row 14) in the "Total"/2022 column
$ws.Range("I14").Value = 673
